# Update 2p0. Convention change to support multi-axle vehicles.
#
# Adds two new vehicle-instance sheets to the Aero coefficients workbook:
#   - "Truck_Amandla"   inserted right after "Bus_Makhulu"
#   - "Trailer_Kumanzi" inserted right after "Trailer_Thwala" (new last sheet)
#
# Both new sheets reuse the same layout/formatting as the existing "sedan"
# class sheet (Bus_Makhulu) -- same styles, column widths, frozen panes,
# tab color -- just with the Instance name (H3) and the sPressureCentre
# vector (F9:H9) updated for the new vehicle instance. Trailer_Kumanzi
# becomes the active/selected sheet, matching the final state on disk.

$wb = $excel.ActiveWorkbook

# --- Insert "Truck_Amandla" right after "Bus_Makhulu" ---------------------
$busSheet = $wb.Worksheets.Item("Bus_Makhulu")
$busSheet.Copy($null, $busSheet)

$truckSheet = $wb.Worksheets.Item("Bus_Makhulu (2)")
$truckSheet.Name = "Truck_Amandla"
$truckSheet.Range("H3").Value = "Truck_Amandla"
$truckSheet.Range("H9").Value = 1.1
$truckSheet.Activate()
$truckSheet.Range("H5:H9").Select()

# --- Insert "Trailer_Kumanzi" right after "Trailer_Thwala" ----------------
$thwalaSheet = $wb.Worksheets.Item("Trailer_Thwala")
$busSheet.Copy($null, $thwalaSheet)

$kumanziSheet = $wb.Worksheets.Item("Bus_Makhulu (2)")
$kumanziSheet.Name = "Trailer_Kumanzi"
$kumanziSheet.Range("H3").Value = "Trailer_Kumanzi"
$kumanziSheet.Range("F9").Value = 5
$kumanziSheet.Range("G9").Value = 0
$kumanziSheet.Range("H9").Value = 2

# Trailer_Kumanzi ends up as the active tab in the saved workbook.
$kumanziSheet.Activate()
$kumanziSheet.Range("J20").Select()
